$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(2)
$rng = $p.Range
Write-Host "Before: [$($rng.Text)]"
$found = $rng.Find.Execute("Video Game Integrative Task 2", $true, $false, $false, $false, $false, $true, 1, $false, "NeoTunes - Integrative Task 3", 2)
Write-Host "Found: $found"
Write-Host "After: [$($d.Paragraphs.Item(2).Range.Text)]"
